# feat: Se agrego tabla de registro de tablas obtenidas de la API
#
# Adds a new "Tablas Extraidas de la API" registry table in rows 21-28 of
# Hoja1, listing the tables pulled from the API and whether they've been
# extracted yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B gets a touch wider to fit the new content ------------------
$ws.Columns("B").ColumnWidth = 13.5546875

# --- Row 21: merged title bar --------------------------------------------
$ws.Range("B21").Value = "Tablas Extraidas de la API"

$titleLeft = $ws.Range("B21")
$titleLeft.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$titleLeft.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$titleLeft.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$titleLeft.HorizontalAlignment = -4108     # xlCenter

$titleRight = $ws.Range("C21")
$titleRight.Borders.Item(10).LineStyle = 1 # xlEdgeRight
$titleRight.Borders.Item(8).LineStyle = 1  # xlEdgeTop
$titleRight.Borders.Item(9).LineStyle = 1  # xlEdgeBottom
$titleRight.HorizontalAlignment = -4108    # xlCenter

$ws.Range("B21:C21").Merge()

# --- Row 22: column headers for the registry table ------------------------
$ws.Range("B5").Copy()
$ws.Range("B22:C22").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B22").Value = "Nombre tabla"
$ws.Range("C22").Value = "Extraida"

# --- Row 23: repurpose the old blank placeholder row -----------------------
$ws.Range("D5").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("B23").Value = "anime"

$c23 = $ws.Range("C23")
$c23.Font.Bold = $true
$c23.Font.Name = "Arial Unicode MS"
$c23.Font.Size = 10
$c23.Borders.LineStyle = 1
$c23.Value = $null

# --- Rows 24-28: remaining table names ------------------------------------
$ws.Range("B5").Copy()
$ws.Range("B24:C28").PasteSpecial(-4122)

$ws.Range("B24").Value = "generos"
$ws.Range("B25").Value = "studios"
$ws.Range("B26").Value = "popularidad"
$ws.Range("B27").Value = "anime_genero"
$ws.Range("B28").Value = "anime_studios"

$ws.Application.CutCopyMode = $false

# --- Row 30: stray formatted cell left over from editing -------------------
$ws.Range("E30").Font.Underline = 2   # xlUnderlineStyleSingle

# --- View state: scrolled down, selection resting below the new table ------
$ws.Range("B31").Select()
$ws.Application.ActiveWindow.ScrollRow = 10
